# Apply the CSAT Performance Reports update:
#  - New feedback data arrived for 2025-10-13 (serial 45943) and 2025-10-14 (serial 45944)
#  - Daywise_Report gets two new day rows (14, 15)
#  - Agentwise_Report gets a new agent "VPS248" inserted in sorted order (row 24),
#    pushing VPS256/VPS257/VPS264 down by one row; existing agent SBM967 (row 13)
#    gets an extra CSAT3 response
#  - Negative_Responses gets one new negative-feedback record (row 13)
#  - Daywise_Agent_Performance gets two new rows (31, 32) for the new responses
#  - All four tables (ListObjects) are resized to cover the new data

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheet "Daywise_Report"  (Table_Daywise_Report)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daywise_Report")

# Update the MTD summary row (row 2) totals
$ws1.Range("B2").Value = 12
$ws1.Range("D2").Value = 16
$ws1.Range("E2").Value = 31
$ws1.Range("F2").Value = 2.13

# New row 14 (2025-10-13): copy formatting from existing same-styled cells, then set values
$ws1.Range("A3").Copy()
$ws1.Range("A14").PasteSpecial($xlPasteFormats)
$ws1.Range("B3").Copy()
$ws1.Range("B14").PasteSpecial($xlPasteFormats)
$ws1.Range("C3").Copy()
$ws1.Range("C14").PasteSpecial($xlPasteFormats)
$ws1.Range("D3").Copy()
$ws1.Range("D14").PasteSpecial($xlPasteFormats)
$ws1.Range("E3").Copy()
$ws1.Range("E14").PasteSpecial($xlPasteFormats)
$ws1.Range("F3").Copy()
$ws1.Range("F14").PasteSpecial($xlPasteFormats)

$ws1.Range("A14").Value = 45943
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 0
$ws1.Range("D14").Value = 1
$ws1.Range("E14").Value = 1
$ws1.Range("F14").Value = 3

# New row 15 (2025-10-14)
$ws1.Range("A3").Copy()
$ws1.Range("A15").PasteSpecial($xlPasteFormats)
$ws1.Range("D3").Copy()
$ws1.Range("B15").PasteSpecial($xlPasteFormats)
$ws1.Range("C3").Copy()
$ws1.Range("C15").PasteSpecial($xlPasteFormats)
$ws1.Range("B3").Copy()
$ws1.Range("D15").PasteSpecial($xlPasteFormats)
$ws1.Range("E3").Copy()
$ws1.Range("E15").PasteSpecial($xlPasteFormats)
$ws1.Range("D3").Copy()
$ws1.Range("F15").PasteSpecial($xlPasteFormats)

$ws1.Range("A15").Value = 45944
$ws1.Range("B15").Value = 1
$ws1.Range("C15").Value = 0
$ws1.Range("D15").Value = 0
$ws1.Range("E15").Value = 1
$ws1.Range("F15").Value = 1

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:F15"))

# ---------------------------------------------------------------------------
# Sheet "Agentwise_Report"  (Table_Agentwise_Report)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Agentwise_Report")

# SBM967 (row 13) got one more CSAT3 response: D13 1->2, change fill red->yellow
$ws1.Range("B6").Copy()
$ws2.Range("D13").PasteSpecial($xlPasteFormats)
$ws2.Range("D13").Value = 2
$ws2.Range("E13").Value = 2

# Insert a new row at position 24 for the new agent "VPS248" (keeps the sheet
# alphabetically sorted); this shifts old rows 24/25/26 (VPS256/VPS257/VPS264)
# down to 25/26/27, preserving their values & formatting.
$ws2.Rows("24:24").Insert()

$ws2.Range("B2").Copy()
$ws2.Range("B24").PasteSpecial($xlPasteFormats)
$ws2.Range("D2").Copy()
$ws2.Range("C24").PasteSpecial($xlPasteFormats)
$ws2.Range("D2").Copy()
$ws2.Range("D24").PasteSpecial($xlPasteFormats)
$ws2.Range("D2").Copy()
$ws2.Range("E24").PasteSpecial($xlPasteFormats)
$ws2.Range("F2").Copy()
$ws2.Range("F24").PasteSpecial($xlPasteFormats)

$ws2.Range("A24").Value = "VPS248"
$ws2.Range("B24").Value = 1
$ws2.Range("C24").Value = 0
$ws2.Range("D24").Value = 0
$ws2.Range("E24").Value = 1
$ws2.Range("F24").Value = 1

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:F27"))

# ---------------------------------------------------------------------------
# Sheet "Negative_Responses"  (Table_Negative_Responses)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Negative_Responses")

$ws3.Range("R2").Copy()
$ws3.Range("R13").PasteSpecial($xlPasteFormats)

$ws3.Range("A13").Value = "J142467211139"
$ws3.Range("B13").Value = 9840814858
$ws3.Range("C13").Value = "14-10-2025 07:49 PM"
$ws3.Range("D13").Value = "CSAT 1"
$ws3.Range("E13").Value = "14-10-2025 08:51 PM"
$ws3.Range("H13").Value = "QUERY"
$ws3.Range("I13").Value = "PREPAID"
$ws3.Range("J13").Value = "SERVICES"
$ws3.Range("K13").Value = "DATA SERVICES"
$ws3.Range("L13").Value = "DATA ACTIVATION"
$ws3.Range("M13").Value = "DHARMAPURI"
$ws3.Range("N13").Value = "TAMILNADU"
$ws3.Range("O13").Value = "TAMILNADU"
$ws3.Range("P13").Value = "VPS248"
$ws3.Range("Q13").Value = 1
$ws3.Range("R13").Value = 45944

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:R13"))

# ---------------------------------------------------------------------------
# Sheet "Daywise_Agent_Performance"  (Table_Daywise_Agent_Performance)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Daywise_Agent_Performance")

# New row 31: SBM967 on 2025-10-13
$ws4.Range("A2").Copy()
$ws4.Range("A31").PasteSpecial($xlPasteFormats)
$ws4.Range("G2").Copy()
$ws4.Range("G31").PasteSpecial($xlPasteFormats)

$ws4.Range("A31").Value = 45943
$ws4.Range("B31").Value = "SBM967"
$ws4.Range("C31").Value = 0
$ws4.Range("D31").Value = 0
$ws4.Range("E31").Value = 1
$ws4.Range("F31").Value = 1
$ws4.Range("G31").Value = 3

# New row 32: VPS248 on 2025-10-14
$ws4.Range("A2").Copy()
$ws4.Range("A32").PasteSpecial($xlPasteFormats)
$ws4.Range("C5").Copy()
$ws4.Range("C32").PasteSpecial($xlPasteFormats)
$ws4.Range("C5").Copy()
$ws4.Range("G32").PasteSpecial($xlPasteFormats)

$ws4.Range("A32").Value = 45944
$ws4.Range("B32").Value = "VPS248"
$ws4.Range("C32").Value = 1
$ws4.Range("D32").Value = 0
$ws4.Range("E32").Value = 0
$ws4.Range("F32").Value = 1
$ws4.Range("G32").Value = 1

$lo4 = $ws4.ListObjects.Item(1)
$lo4.Resize($ws4.Range("A1:G32"))
